# Weekly update: insert a new data row for "Apio" (Femacal de La Calera)
# at row 200, pushing all subsequent rows down by one (old row 200 becomes
# new row 201, ..., old row 307 becomes new row 308).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 200; existing rows 200:307 shift to 201:308.
$ws.Rows("200:200").Insert()

# Populate the newly inserted row 200 with this week's data.
$ws.Range("A200").Value = 3
$ws.Range("B200").Value = "Femacal de La Calera"
$ws.Range("C200").Value = "Coquimbo"
$ws.Range("D200").Value = 44572
$ws.Range("E200").Value = 5
$ws.Range("F200").Value = 100112017
$ws.Range("G200").Value = "Apio"
$ws.Range("H200").Value = "Americana (o)"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 200
$ws.Range("K200").Value = 9000
$ws.Range("L200").Value = 9500
$ws.Range("M200").Value = 9200
$ws.Range("N200").Value = "$/docena de matas"
$ws.Range("O200").Value = "Pan de Azúcar"
$ws.Range("P200").Value = 1533
$ws.Range("Q200").Value = 6
$ws.Range("R200").Value = "Hortaliza"
